# #5: property boat&car done
# Fill out the 汽車 (car) sheet's header row with proper field names and
# populate the new "capacity" column + the property/legislator metadata
# columns (H:N), matching the pattern already used on the other sheets
# (e.g. 股票/stock).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 汽車

# ---- Row 1: header labels ----
$ws.Cells.Item(1,2).Value  = "name"
$ws.Cells.Item(1,3).Value  = "capacity"
$ws.Cells.Item(1,4).Value  = "owner"
$ws.Cells.Item(1,5).Value  = "register_date"
$ws.Cells.Item(1,6).Value  = "register_reason"
$ws.Cells.Item(1,7).Value  = "acquire_value"
$ws.Cells.Item(1,8).Value  = "property_category"
$ws.Cells.Item(1,9).Value  = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

# ---- Row 2 (record index 54) ----
$ws.Cells.Item(2,2).Value  = "BMW740LISEDAN"
$ws.Cells.Item(2,3).Value  = 4000
$ws.Cells.Item(2,4).Value  = "李鴻鈞"
$ws.Cells.Item(2,5).Value  = "劉素幸"
$ws.Cells.Item(2,6).Value  = "買賣"
$ws.Cells.Item(2,7).Value  = 4200000
$ws.Cells.Item(2,8).Value  = "land"
$ws.Cells.Item(2,9).Value  = "normal"
$ws.Cells.Item(2,10).Value = "2012-04-23"
$ws.Cells.Item(2,11).Value = "李鴻鈞"
$ws.Cells.Item(2,12).Value = 898
$ws.Cells.Item(2,13).Value = "tmp651e1"
$ws.Cells.Item(2,14).Value = 54

# ---- Row 3 (record index 55) ----
$ws.Cells.Item(3,2).Value  = "96年02月08日"
$ws.Cells.Item(3,3).Value  = 3456
$ws.Cells.Item(3,4).Value  = "LEXUSES350"
$ws.Cells.Item(3,5).Value  = "97年04月24H"
$ws.Cells.Item(3,6).Value  = "買賣"
$ws.Cells.Item(3,7).Value  = 2050000
$ws.Cells.Item(3,8).Value  = "land"
$ws.Cells.Item(3,9).Value  = "normal"
$ws.Cells.Item(3,10).Value = "2012-04-23"
$ws.Cells.Item(3,11).Value = "李鴻鈞"
$ws.Cells.Item(3,12).Value = 898
$ws.Cells.Item(3,13).Value = "tmp651e1"
$ws.Cells.Item(3,14).Value = 55

Write-Output "car sheet updated"
